# Update values in Sheet1 per the diff (columns B and C for several rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 5.734999999999999
$ws.Range("C4").Value = -12.74

$ws.Range("C5").Value = -12.927

$ws.Range("B6").Value = 6.164
$ws.Range("C6").Value = -12.284

$ws.Range("B7").Value = 6.638

$ws.Range("B8").Value = 5.323
$ws.Range("C8").Value = -12.597

$ws.Range("B16").Value = 5.322
$ws.Range("C16").Value = -12.597

$ws.Range("B20").Value = 5.766999999999999

$ws.Range("B21").Value = 6.247

$ws.Range("C22").Value = -12.78
